# Generate Report for Handback
# Updates the "Latest HO Xliff Generate Date", "Priority", "Correspond Handoff
# Datetime" and "Correspond Handback DateTime" values to reflect a freshly
# regenerated handback status report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet - column G: Latest HO Xliff Generate Date
$wsOverview.Range("G2").Value = "2016-09-05 14:20:47"
$wsOverview.Range("G3").Value = "2016-09-05 14:20:47"

# zh-cn sheet - column E: Priority, column H: Correspond Handoff Datetime,
# column K: Correspond Handback DateTime
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-09-05 14:20:42"
$wsZhCn.Range("H3").Value = "2016-09-05 14:20:42"
$wsZhCn.Range("K2").Value = "2016-09-05 14:21:01"
$wsZhCn.Range("K3").Value = "2016-09-05 14:21:01"

# de-de sheet - column E: Priority, column H: Correspond Handoff Datetime,
# column K: Correspond Handback DateTime
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H2").Value = "2016-09-05 14:20:47"
$wsDeDe.Range("H3").Value = "2016-09-05 14:20:47"
$wsDeDe.Range("K2").Value = "2016-09-05 14:21:18"
$wsDeDe.Range("K3").Value = "2016-09-05 14:21:18"
